$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update word list / image list / category labels for rows 2-49 (row 1 header unchanged)
$ws.Range("A2").Value = 'wenden'
$ws.Range("B2").Value = 'none'
$ws.Range("C2").Value = 'none'
$ws.Range("A3").Value = 'handeln'
$ws.Range("B3").Value = 'house/house012.jpg'
$ws.Range("C3").Value = 'house'
$ws.Range("A4").Value = 'retten'
$ws.Range("B4").Value = 'house/house021.jpg'
$ws.Range("C4").Value = 'house'
$ws.Range("A5").Value = 'rufen'
$ws.Range("B5").Value = 'none'
$ws.Range("C5").Value = 'none'
$ws.Range("A6").Value = 'mauern'
$ws.Range("B6").Value = 'dog/dog009.jpg'
$ws.Range("C6").Value = 'dog'
$ws.Range("A7").Value = 'streichen'
$ws.Range("B7").Value = 'dog/dog008.jpg'
$ws.Range("C7").Value = 'dog'
$ws.Range("A8").Value = 'enden'
$ws.Range("B8").Value = 'none'
$ws.Range("C8").Value = 'none'
$ws.Range("A9").Value = 'lächeln'
$ws.Range("B9").Value = 'house/house009.jpg'
$ws.Range("C9").Value = 'house'
$ws.Range("A10").Value = 'wohnen'
$ws.Range("B10").Value = 'dog/dog007.jpg'
$ws.Range("C10").Value = 'dog'
$ws.Range("A11").Value = 'opfern'
$ws.Range("B11").Value = 'none'
$ws.Range("C11").Value = 'none'
$ws.Range("A12").Value = 'fordern'
$ws.Range("B12").Value = 'house/house000.jpg'
$ws.Range("C12").Value = 'house'
$ws.Range("A13").Value = 'decken'
$ws.Range("B13").Value = 'dog/dog027.jpg'
$ws.Range("C13").Value = 'dog'
$ws.Range("A14").Value = 'schalten'
$ws.Range("B14").Value = 'none'
$ws.Range("C14").Value = 'none'
$ws.Range("A15").Value = 'helfen'
$ws.Range("B15").Value = 'dog/dog022.jpg'
$ws.Range("C15").Value = 'dog'
$ws.Range("A16").Value = 'tanzen'
$ws.Range("B16").Value = 'dog/dog005.jpg'
$ws.Range("C16").Value = 'dog'
$ws.Range("A17").Value = 'klagen'
$ws.Range("B17").Value = 'none'
$ws.Range("C17").Value = 'none'
$ws.Range("A18").Value = 'lehnen'
$ws.Range("B18").Value = 'dog/dog020.jpg'
$ws.Range("C18").Value = 'dog'
$ws.Range("A19").Value = 'schweben'
$ws.Range("B19").Value = 'dog/dog012.jpg'
$ws.Range("C19").Value = 'dog'
$ws.Range("A20").Value = 'sparen'
$ws.Range("B20").Value = 'none'
$ws.Range("C20").Value = 'none'
$ws.Range("A21").Value = 'hoffen'
$ws.Range("B21").Value = 'dog/dog014.jpg'
$ws.Range("C21").Value = 'dog'
$ws.Range("A22").Value = 'zahlen'
$ws.Range("B22").Value = 'house/house026.jpg'
$ws.Range("C22").Value = 'house'
$ws.Range("A23").Value = 'hören'
$ws.Range("B23").Value = 'none'
$ws.Range("C23").Value = 'none'
$ws.Range("A24").Value = 'posten'
$ws.Range("B24").Value = 'house/house011.jpg'
$ws.Range("C24").Value = 'house'
$ws.Range("A25").Value = 'leiten'
$ws.Range("B25").Value = 'dog/dog003.jpg'
$ws.Range("C25").Value = 'dog'
$ws.Range("A26").Value = 'drohen'
$ws.Range("B26").Value = 'none'
$ws.Range("C26").Value = 'none'
$ws.Range("A27").Value = 'bitten'
$ws.Range("B27").Value = 'dog/dog015.jpg'
$ws.Range("C27").Value = 'dog'
$ws.Range("A28").Value = 'kümmern'
$ws.Range("B28").Value = 'dog/dog028.jpg'
$ws.Range("C28").Value = 'dog'
$ws.Range("A29").Value = 'orten'
$ws.Range("B29").Value = 'none'
$ws.Range("C29").Value = 'none'
$ws.Range("A30").Value = 'klingen'
$ws.Range("B30").Value = 'house/house014.jpg'
$ws.Range("C30").Value = 'house'
$ws.Range("A31").Value = 'kleben'
$ws.Range("B31").Value = 'house/house016.jpg'
$ws.Range("C31").Value = 'house'
$ws.Range("A32").Value = 'weigern'
$ws.Range("B32").Value = 'none'
$ws.Range("C32").Value = 'none'
$ws.Range("A33").Value = 'schreiben'
$ws.Range("B33").Value = 'house/house013.jpg'
$ws.Range("C33").Value = 'house'
$ws.Range("A34").Value = 'schneiden'
$ws.Range("B34").Value = 'house/house023.jpg'
$ws.Range("C34").Value = 'house'
$ws.Range("A35").Value = 'ächzen'
$ws.Range("B35").Value = 'none'
$ws.Range("C35").Value = 'none'
$ws.Range("A36").Value = 'leugnen'
$ws.Range("B36").Value = 'house/house019.jpg'
$ws.Range("C36").Value = 'house'
$ws.Range("A37").Value = 'deuten'
$ws.Range("B37").Value = 'dog/dog025.jpg'
$ws.Range("C37").Value = 'dog'
$ws.Range("A38").Value = 'kosten'
$ws.Range("B38").Value = 'none'
$ws.Range("C38").Value = 'none'
$ws.Range("A39").Value = 'spenden'
$ws.Range("B39").Value = 'dog/dog006.jpg'
$ws.Range("C39").Value = 'dog'
$ws.Range("A40").Value = 'leisten'
$ws.Range("B40").Value = 'house/house025.jpg'
$ws.Range("C40").Value = 'house'
$ws.Range("A41").Value = 'dauern'
$ws.Range("B41").Value = 'none'
$ws.Range("C41").Value = 'none'
$ws.Range("A42").Value = 'stürmen'
$ws.Range("B42").Value = 'dog/dog001.jpg'
$ws.Range("C42").Value = 'dog'
$ws.Range("A43").Value = 'fühlen'
$ws.Range("B43").Value = 'house/house020.jpg'
$ws.Range("C43").Value = 'house'
$ws.Range("A44").Value = 'stören'
$ws.Range("B44").Value = 'none'
$ws.Range("C44").Value = 'none'
$ws.Range("A45").Value = 'binden'
$ws.Range("B45").Value = 'house/house003.jpg'
$ws.Range("C45").Value = 'house'
$ws.Range("A46").Value = 'gelten'
$ws.Range("B46").Value = 'house/house029.jpg'
$ws.Range("C46").Value = 'house'
$ws.Range("A47").Value = 'bremsen'
$ws.Range("B47").Value = 'none'
$ws.Range("C47").Value = 'none'
$ws.Range("A48").Value = 'geben'
$ws.Range("B48").Value = 'house/house024.jpg'
$ws.Range("C48").Value = 'house'
$ws.Range("A49").Value = 'holen'
$ws.Range("B49").Value = 'dog/dog031.jpg'
$ws.Range("C49").Value = 'dog'
